$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# New header for the added `test_xlr_n_percent` column (K)
$ws.Range("K2").Value = "test_xlr_n_percent"

# xlRight / xlBottom alignment constants
$xlRight = -4152
$xlBottom = -4107

# Build the style on the first data row, then stamp (copy) that exact
# style down the rest of the column so the cellXfs table doesn't grow
# one entry per cell (font/alignment get de-duplicated against the
# first cell's new style).
$first = $ws.Cells.Item(3, 11)
$first.Value = "1 (3%)"
$first.HorizontalAlignment = $xlRight
$first.VerticalAlignment = $xlBottom
$first.Font.Name = "calibri"

$first.Copy()
$ws.Range("K4:K35").PasteSpecial(-4122)  # xlPasteFormats

$values = @(
  "2 (6%)","3 (9%)","4 (12%)","5 (16%)","6 (19%)","7 (22%)","8 (25%)","9 (28%)","10 (31%)",
  "11 (34%)","12 (38%)","13 (41%)","14 (44%)","15 (47%)","16 (50%)","17 (53%)","18 (56%)","19 (59%)","20 (62%)",
  "21 (66%)","22 (69%)","23 (72%)","24 (75%)","25 (78%)","26 (81%)","27 (84%)","28 (88%)","29 (91%)","30 (94%)",
  "31 (97%)","32 (100%)"
)
for ($i = 0; $i -lt $values.Length; $i++) {
    $row = $i + 4
    $ws.Cells.Item($row, 11).Value = $values[$i]
}

# Row 35 keeps the style but has no value, matching the rest of the
# (otherwise-empty) trailing row.
$ws.Cells.Item(35, 11).Value = ""
